$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.085.49'
$ws.Range("E2").Value = '  -1.24%  '
$ws.Range("D3").Value = '3.516.80'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '585.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.61'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.78%  '
$ws.Range("D7").Value = '3.516.81'
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("E10").Value = '  -0.40%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.11'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.63%  '
$ws.Range("E12").Value = '  -2.18%  '
$ws.Range("D13").Value = '4.117.56'
$ws.Range("E13").Value = '  +0.10%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.45'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.33%  '
$ws.Range("E15").Value = '  +1.34%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.524.65'
$ws.Range("E16").Value = '  +0.22%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000179'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.05%  '
$ws.Range("D18").Value = '64.143.97'
$ws.Range("E18").Value = '  -1.20%  '
$ws.Range("E19").Value = '  -2.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.86'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.83%  '
$ws.Range("E21").Value = '  -0.52%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '382.95'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.87%  '
$ws.Range("D23").Value = '3.660.73'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.568'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.93%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.10'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.99%  '
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.71'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.54%  '
$ws.Range("E28").Value = '  +3.71%  '
$ws.Range("E29").Value = '  -1.30%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.48'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.03%  '
$ws.Range("E32").Value = '  +1.50%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.22'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.29%  '
$ws.Range("D34").Value = '3.533.49'
$ws.Range("E34").Value = '  +0.36%  '
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '23.54'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.19%  '
$ws.Range("E37").Value = '  -0.62%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.41'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.63%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.93'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.25%  '
$ws.Range("E40").Value = '  -0.61%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '160.13'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.66%  '
$ws.Range("E42").Value = '  -2.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.60'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.04%  '
$ws.Range("E44").Value = '  -0.82%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '41.62'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.95%  '
$ws.Range("E47").Value = '  -3.12%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.39'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.30%  '
$ws.Range("E49").Value = '  -2.51%  '
$ws.Range("D50").Value = '2.482.23'
$ws.Range("E50").Value = '  -0.15%  '
$ws.Range("E51").Value = '  -0.78%  '
